$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NIG(0.9055176893981262, 0.5987566529066475, 1.337029383951728, 3.065098492186134)"
$ws.Range("C2").Value = "JSU(-1.1417515853983575, 1.0599534581957901, 2.742539681853493, 3.604525022500561)"
$ws.Range("D2").Value = "NIG(0.9065012229598297, 0.6863998286522617, 1.242806036771225, 3.2838978421660214)"
$ws.Range("E2").Value = "JSU(-0.9294905653722012, 1.2603161415622215, 4.425964563092953, 6.201963262063888)"
